$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 6) of trading/sentiment data, mirroring the
# existing rows 2-5 in xl/worksheets/sheet1.xml.
$ws.Range("A6").Value = 42606.880983796298
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 42
$ws.Range("E6").Value = 61
$ws.Range("F6").Value = 38
$ws.Range("G6").Value = 4309
$ws.Range("H6").Value = 10396
$ws.Range("I6").Value = 1203
$ws.Range("J6").Value = 115
$ws.Range("K6").Value = 87
$ws.Range("L6").Value = 8
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = "Bag"
